$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date (row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-02 07:14:22"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-02 07:14:17"
$wsZhCn.Range("K2").Value = "2016-09-02 07:14:34"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime (row 2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-02 07:14:22"
$wsDeDe.Range("K2").Value = "2016-09-02 07:14:41"
